# Daily attendance processing - 2026-01-05 23:02:16
# Reverses the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $rev = $trimmed[($trimmed.Count - 1)..0]
        $newVal = [string]::Join(", ", $rev)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
